$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 144, shifting rows 144:233 down to 145:234
$ws.Rows(144).Insert()

$ws.Range("A144").Value = 10
$ws.Range("B144").Value = 'Vega Modelo de Temuco'
$ws.Range("C144").Value = 'La Araucanía'
$ws.Range("D144").Value = 45176
$ws.Range("E144").Value = 9
$ws.Range("F144").Value = 100114002
$ws.Range("G144").Value = 'Camote'
$ws.Range("H144").Value = 'Sin especificar'
$ws.Range("I144").Value = 'Primera'
$ws.Range("J144").Value = 80
$ws.Range("K144").Value = 24000
$ws.Range("L144").Value = 24000
$ws.Range("M144").Value = 24000
$ws.Range("N144").Value = '$/caja 18 kilos'
$ws.Range("O144").Value = 'Perú'
$ws.Range("P144").Value = 1333
$ws.Range("Q144").Value = 18
$ws.Range("R144").Value = 'Hortaliza'
